# Auto-generated Excel COM-interop edit script
# Applies targeted cell value updates across multiple worksheets
# to match the committed timetable / classroom allocation changes.

$wb = $excel.ActiveWorkbook

# ---- Sheet: Regular_Timetable ----
$ws = $wb.Worksheets.Item("Regular_Timetable")
$ws.Range("D23").Value = "Tue 09:00-10:30, Thu 13:00-14:30 [L406]"
$ws.Range("E23").Value = "Wed 14:30-15:30 [L408]"
$ws.Range("D24").Value = "Tue 09:00-10:30, Thu 13:00-14:30 [L408]"
$ws.Range("E24").Value = "Wed 14:30-15:30"
$ws.Range("D25").Value = "Tue 09:00-10:30, Thu 13:00-14:30"
$ws.Range("E25").Value = "Wed 14:30-15:30"
$ws.Range("D26").Value = "Tue 09:00-10:30, Thu 13:00-14:30"
$ws.Range("E26").Value = "Wed 14:30-15:30"
$ws.Range("E29").Value = "Thu 14:30-15:30"
$ws.Range("D32").Value = "Tue 15:30-17:00 [C302], Thu 15:30-17:00 [C302]"
$ws.Range("D33").Value = "Tue 15:30-17:00 [C304], Thu 15:30-17:00 [C304]"
$ws.Range("D34").Value = "Tue 15:30-17:00 [C305], Thu 15:30-17:00 [C305]"
$ws.Range("D35").Value = "Tue 15:30-17:00 [L402], Thu 15:30-17:00 [L402]"

# ---- Sheet: PreMid_Timetable ----
$ws = $wb.Worksheets.Item("PreMid_Timetable")
$ws.Range("D23").Value = "Tue 09:00-10:30, Thu 13:00-14:30 [L406]"
$ws.Range("E23").Value = "Wed 14:30-15:30 [L408]"
$ws.Range("D24").Value = "Tue 09:00-10:30, Thu 13:00-14:30 [L408]"
$ws.Range("E24").Value = "Wed 14:30-15:30"
$ws.Range("D25").Value = "Tue 09:00-10:30, Thu 13:00-14:30"
$ws.Range("E25").Value = "Wed 14:30-15:30"
$ws.Range("D26").Value = "Tue 09:00-10:30, Thu 13:00-14:30"
$ws.Range("E26").Value = "Wed 14:30-15:30"
$ws.Range("E29").Value = "Thu 14:30-15:30"
$ws.Range("D32").Value = "Tue 15:30-17:00 [C302], Thu 15:30-17:00 [C302]"
$ws.Range("D33").Value = "Tue 15:30-17:00 [C304], Thu 15:30-17:00 [C304]"
$ws.Range("D34").Value = "Tue 15:30-17:00 [C305], Thu 15:30-17:00 [C305]"
$ws.Range("D35").Value = "Tue 15:30-17:00 [L402], Thu 15:30-17:00 [L402]"

# ---- Sheet: PostMid_Timetable ----
$ws = $wb.Worksheets.Item("PostMid_Timetable")
$ws.Range("D23").Value = "Tue 09:00-10:30, Thu 13:00-14:30 [L406]"
$ws.Range("E23").Value = "Wed 14:30-15:30 [L408]"
$ws.Range("D24").Value = "Tue 09:00-10:30, Thu 13:00-14:30 [L408]"
$ws.Range("E24").Value = "Wed 14:30-15:30"
$ws.Range("D25").Value = "Tue 09:00-10:30, Thu 13:00-14:30"
$ws.Range("E25").Value = "Wed 14:30-15:30"
$ws.Range("D26").Value = "Tue 09:00-10:30, Thu 13:00-14:30"
$ws.Range("E26").Value = "Wed 14:30-15:30"
$ws.Range("E29").Value = "Thu 14:30-15:30"
$ws.Range("D32").Value = "Tue 15:30-17:00 [C302], Thu 15:30-17:00 [C302]"
$ws.Range("D33").Value = "Tue 15:30-17:00 [C304], Thu 15:30-17:00 [C304]"
$ws.Range("D34").Value = "Tue 15:30-17:00 [C305], Thu 15:30-17:00 [C305]"
$ws.Range("D35").Value = "Tue 15:30-17:00 [L402], Thu 15:30-17:00 [L402]"

# ---- Sheet: Section_A ----
$ws = $wb.Worksheets.Item("Section_A")
$ws.Range("B16").Value = "Mini Project [C102]"
$ws.Range("C16").Value = "0-0-0-8-2 [C102]"
$ws.Range("D16").Value = "Full Sem [C102]"
$ws.Range("E16").Value = "0/0 [C102]"
$ws.Range("F16").Value = "0/0 [C102]"
$ws.Range("D32").Value = "Tue 15:30-17:00 [C302], Thu 15:30-17:00"
$ws.Range("D33").Value = "Tue 15:30-17:00 [C302], Thu 15:30-17:00"
$ws.Range("D34").Value = "Tue 15:30-17:00 [C302], Thu 15:30-17:00"
$ws.Range("D35").Value = "Tue 15:30-17:00 [C302], Thu 15:30-17:00"

# ---- Sheet: Classroom_Utilization ----
$ws = $wb.Worksheets.Item("Classroom_Utilization")
$ws.Range("D7").Value = 7.5
$ws.Range("E7").Value = 1.5
$ws.Range("D23").Value = 6
$ws.Range("E23").Value = 1.2
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0

# ---- Sheet: Classroom_Allocation ----
$ws = $wb.Worksheets.Item("Classroom_Allocation")
$ws.Range("I2").Value = "Projector"
$ws.Range("M2").Value = "C102"
$ws.Range("I3").Value = "Projector"
$ws.Range("M3").Value = "C102"
$ws.Range("I4").Value = "Projector"
$ws.Range("M4").Value = "C102"
$ws.Range("I5").Value = "Projector"
$ws.Range("M5").Value = "C102"
$ws.Range("I6").Value = "Projector"
$ws.Range("M6").Value = "C102"
$ws.Range("I7").Value = "Projector"
$ws.Range("M7").Value = "C101"
$ws.Range("G8").Value = "classroom"
$ws.Range("H8").Value = "96"
$ws.Range("I8").Value = "Projector"
$ws.Range("M8").Value = "C102"
$ws.Range("G9").Value = "classroom"
$ws.Range("H9").Value = "96"
$ws.Range("M9").Value = "C104"
$ws.Range("H15").Value = "120"
$ws.Range("M15").Value = "C002"
$ws.Range("H16").Value = "120"
$ws.Range("M16").Value = "C002"
$ws.Range("G17").Value = "Auditorium"
$ws.Range("H17").Value = "240"
$ws.Range("I17").Value = "Audio/Video System"
$ws.Range("M17").Value = "C004"
$ws.Range("M18").Value = "L406"
$ws.Range("M20").Value = "L407"
$ws.Range("H32").Value = "120"
$ws.Range("M32").Value = "C002"
$ws.Range("H33").Value = "120"
$ws.Range("M33").Value = "C002"
$ws.Range("G34").Value = "Auditorium"
$ws.Range("H34").Value = "240"
$ws.Range("I34").Value = "Audio/Video System"
$ws.Range("M34").Value = "C004"
$ws.Range("M35").Value = "L406"
$ws.Range("M37").Value = "L406"
$ws.Range("M46").Value = "L408"
$ws.Range("H47").Value = "96"
$ws.Range("I47").Value = "Projector"
$ws.Range("M47").Value = "C101"
$ws.Range("H48").Value = "96"
$ws.Range("I48").Value = "Projector"
$ws.Range("M48").Value = "C102"
$ws.Range("H49").Value = "96"
$ws.Range("I49").Value = "Projector"
$ws.Range("M49").Value = "C104"
$ws.Range("H50").Value = "120"
$ws.Range("M50").Value = "C002"
$ws.Range("H51").Value = "120"
$ws.Range("M51").Value = "C002"
$ws.Range("G52").Value = "Auditorium"
$ws.Range("H52").Value = "240"
$ws.Range("I52").Value = "Audio/Video System"
$ws.Range("M52").Value = "C004"
$ws.Range("M53").Value = "L405"
$ws.Range("F55").Value = "Tue 15:30-17:00 [C302], Thu 15:30-17:00"
$ws.Range("M55").Value = "L406"
$ws.Range("N55").Value = "Tue 15:30-17:00 [C302], Thu 15:30-17:00"
$ws.Range("M56").Value = "L406"
$ws.Range("M57").Value = "L408"
$ws.Range("H58").Value = "96"
$ws.Range("I58").Value = "Projector"
$ws.Range("M58").Value = "C101"
$ws.Range("H59").Value = "96"
$ws.Range("I59").Value = "Projector"
$ws.Range("M59").Value = "C102"
$ws.Range("H62").Value = "96"
$ws.Range("I62").Value = "Projector"
$ws.Range("M62").Value = "C104"
$ws.Range("M63").Value = "C202"
$ws.Range("I64").Value = "TV"
$ws.Range("M64").Value = "C203"
$ws.Range("M65").Value = "C302"
$ws.Range("M66").Value = "C304"
$ws.Range("I67").Value = "TV"
$ws.Range("M67").Value = "C305"
$ws.Range("H68").Value = "80"
$ws.Range("M68").Value = "L402"
$ws.Range("H69").Value = "120"
$ws.Range("M69").Value = "C002"
$ws.Range("H70").Value = "120"
$ws.Range("I70").Value = ""
$ws.Range("M70").Value = "C001"
$ws.Range("G71").Value = "Auditorium"
$ws.Range("H71").Value = "240"
$ws.Range("I71").Value = "Audio/Video System"
$ws.Range("M71").Value = "C004"
$ws.Range("M72").Value = "L405"
$ws.Range("M74").Value = "L406"
$ws.Range("I75").Value = "Projector"
$ws.Range("M75").Value = "C101"
$ws.Range("G76").Value = "classroom"
$ws.Range("H76").Value = "96"
$ws.Range("I76").Value = "Projector"
$ws.Range("M76").Value = "C102"
$ws.Range("G77").Value = "classroom"
$ws.Range("H77").Value = "96"
$ws.Range("M77").Value = "C104"
$ws.Range("G78").Value = "classroom"
$ws.Range("H78").Value = "96"
$ws.Range("M78").Value = "C202"
$ws.Range("H79").Value = "120"
$ws.Range("M79").Value = "C002"
$ws.Range("H80").Value = "120"
$ws.Range("M80").Value = "C002"
$ws.Range("G81").Value = "Auditorium"
$ws.Range("H81").Value = "240"
$ws.Range("I81").Value = "Audio/Video System"
$ws.Range("M81").Value = "C004"

# ---- Sheet: Basket_Course_Allocations ----
$ws = $wb.Worksheets.Item("Basket_Course_Allocations")
$ws.Range("C2").Value = "C004, C101"
$ws.Range("C3").Value = "C102"
$ws.Range("C4").Value = "C104"
$ws.Range("C7").Value = "C104"
$ws.Range("C8").Value = "C202"
$ws.Range("C9").Value = "C203"
$ws.Range("C10").Value = "C101, L406, L408"
$ws.Range("C11").Value = "C101, C102, L408"
$ws.Range("C12").Value = "C101, C102, C104"
$ws.Range("C13").Value = "C102, C104, C202"
$ws.Range("C14").Value = "C101, C302, L406"
$ws.Range("C15").Value = "C102, C304"
$ws.Range("C16").Value = "C104, C305"
$ws.Range("C17").Value = "C202, L402"

